# Add a new column K ("tm_vettor") to the ESP_FATT sheet:
#  - K1 header "tm_vettor" using the same style as the other header cells (J1)
#  - K2:K101 filled with the text value "3" using the same style as the other
#    data cells (J2:J101)
#  - column K given the same width as the other (customWidth) columns
#  - dimension grows from A1:J101 to A1:K101 (handled automatically by Excel)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- values first (so the later format paste does not get clobbered) ---

# Header cell
$ws.Range("K1").Value = "tm_vettor"

# Data cells: write as text "3" (quote-prefix forces text, matching the
# existing text-typed sibling columns instead of a numeric 3)
$ws.Range("K2:K101").Value = "'3"

# --- formatting: copy from column J, which already carries the header /
# data styles we want (s=2 for row 1, s=5 for the rest) ---
$ws.Range("J1:J101").Copy() | Out-Null
$ws.Range("K1:K101").PasteSpecial(-4122) | Out-Null

# --- column width: match the other (custom-width) columns ---
$ws.Columns.Item(11).ColumnWidth = 13.08

Write-Output "done"
